# Update the cryptos list (prices/volumes and some coin re-ordering) per GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '89.420.22'
$ws.Range("E2").Value = '  -2.06%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.070.89'
$ws.Range("E3").Value = '  -3.60%  '
$ws.Range("E4").Value = '  +0.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.78'
$ws.Range("E5").Value = '  -1.95%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '613.12'
$ws.Range("E6").Value = '  -3.47%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.366'
$ws.Range("E7").Value = '  -11.07%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.874'
$ws.Range("E8").Value = '  +19.64%  '
$ws.Range("E9").Value = '  +0.19%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '3.070.54'
$ws.Range("E10").Value = '  -3.54%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.661'
$ws.Range("E11").Value = '  +16.27%  '
$ws.Range("E12").Value = '  +2.07%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000239'
$ws.Range("E13").Value = '  -8.07%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.35'
$ws.Range("E14").Value = '  +0.29%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '89.487.45'
$ws.Range("E15").Value = '  -1.60%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '32.27'
$ws.Range("E16").Value = '  -1.61%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.644.38'
$ws.Range("E17").Value = '  -3.28%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.127.43'
$ws.Range("E18").Value = '  -1.77%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.39'
$ws.Range("E19").Value = '  +1.85%  '
$ws.Range("E20").Value = '  -5.67%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.51'
$ws.Range("E21").Value = '  +0.58%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '429.77'
$ws.Range("E22").Value = '  -0.85%  '
$ws.Range("B23").Value = 'Polkadot'
$ws.Range("C23").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.00'
$ws.Range("E23").Value = '  +0.03%  '
$ws.Range("B24").Value = 'Uniswap'
$ws.Range("C24").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '8.17'
$ws.Range("E24").Value = '  -3.88%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '5.46'
$ws.Range("E25").Value = '  +2.36%  '
$ws.Range("B26").Value = 'Aptos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.99'
$ws.Range("E26").Value = '  +2.03%  '
$ws.Range("B27").Value = 'Litecoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '84.78'
$ws.Range("E27").Value = '  +4.83%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.244.95'
$ws.Range("E28").Value = '  -3.19%  '
$ws.Range("E29").Value = '  +0.20%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.164'
$ws.Range("E31").Value = '  +3.73%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.14'
$ws.Range("E32").Value = '  -2.96%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '507.69'
$ws.Range("E33").Value = '  -1.61%  '
$ws.Range("E34").Value = '  -12.20%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.64'
$ws.Range("E35").Value = '  -5.05%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '22.70'
$ws.Range("E36").Value = '  +1.36%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.24'
$ws.Range("E37").Value = '  -3.41%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.78'
$ws.Range("E38").Value = '  -5.74%  '
$ws.Range("E39").Value = '  +3.23%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '22.28'
$ws.Range("E40").Value = '  -0.31%  '
$ws.Range("E41").Value = '  +0.27%  '
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("B43").Value = 'PolygonEcosystemToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.369'
$ws.Range("E43").Value = '  -0.95%  '
$ws.Range("B44").Value = 'Stellar'
$ws.Range("C44").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.138'
$ws.Range("E44").Value = '  +9.93%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.83'
$ws.Range("E45").Value = '  -5.02%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '147.24'
$ws.Range("E46").Value = '  -0.08%  '
$ws.Range("B47").Value = 'Hedera'
$ws.Range("C47").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0697'
$ws.Range("E47").Value = '  +13.04%  '
$ws.Range("B48").Value = 'OKB'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '43.69'
$ws.Range("E48").Value = '  -0.81%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.06'
$ws.Range("E49").Value = '  +0.85%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.21'
$ws.Range("E50").Value = '  +0.25%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '156.30'
$ws.Range("E51").Value = '  -7.83%  '
